$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Change shared string "test" -> "Test1" (group column header value used by D2:D5)
$ws.Range("D2:D5").Value = "Test1"

# 2. Update token column (C2:C5) values: subtract 10
$ws.Range("C2").Value = 3
$ws.Range("C3").Value = 4
$ws.Range("C4").Value = 5
$ws.Range("C5").Value = 6

# 3. Apply new font style (Arial 10, black) to D2:D5
$ws.Range("D2:D5").Font.Name = "Arial"
$ws.Range("D2:D5").Font.Size = 10
$ws.Range("D2:D5").Font.Color = 0

# 4. Update selection to C5
$ws.Range("C5").Select()

# 5. Set page orientation to portrait
$ws.PageSetup.Orientation = 1
